# OLX Monitor update - 2026-02-18 09:33
# The monitoring script appends one new snapshot block every run. This run's
# block (rows 68-74) duplicates the previous block (rows 61-67) verbatim
# (same listings, same formatting) except for the refreshed check timestamp
# in column A and the "days listed" counter in column F, which ticks up for
# listings that are still online.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the last snapshot block (values + styles) into the new rows.
$src = $ws.Range("A61:H67")
$dst = $ws.Range("A68:H74")
$src.Copy($dst)

# Refresh the "last checked" timestamp for every row in the new block.
$timestamp = "2026-02-18 09:33:13"
$ws.Range("A68:A74").Value = $timestamp

# "Przytulny pokój blisko Politechniki" has been online one more day.
$ws.Cells.Item(70, 6).Value = 131
